$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-26 down to 25-27.
$ws.Rows.Item(24).Insert()

# Copy the date number format (style) from the cell that used to be row 24 (now row 25)
# onto the newly inserted row 24's date cell.
$ws.Range("D25").Copy()
$ws.Range("D24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Populate the new row 24 with the new weekly data point.
$ws.Cells.Item(24, 1).Value = 1
$ws.Cells.Item(24, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(24, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(24, 4).Value = 44769
$ws.Cells.Item(24, 5).Value = 15
$ws.Cells.Item(24, 6).Value = 100112043
$ws.Cells.Item(24, 7).Value = "Pepino dulce"
$ws.Cells.Item(24, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 140
$ws.Cells.Item(24, 11).Value = 17000
$ws.Cells.Item(24, 12).Value = 18000
$ws.Cells.Item(24, 13).Value = 17500
$ws.Cells.Item(24, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 16).Value = 972
$ws.Cells.Item(24, 17).Value = 18
$ws.Cells.Item(24, 18).Value = "Hortaliza"

$wb.Save()
